# "added mentor checkbox funcitonality"
# Append a new confirmed-match row: mentor "Henkel" paired with solver
# "Elpis Solar" (COMMENTS left blank for this row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIT_SOLVE_Confirmed_Matches")

$ws.Range("A2").Value = "Henkel"
$ws.Range("B2").Value = "Elpis Solar"
